$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update TestCases (B2) from 45 -> 48
$ws.Range("B2").Value = "48"

# Update Instance (D2) from Automation2 -> Automation3
$ws.Range("D2").Value = "Automation3"

# Move active cell selection to D2
$ws.Range("D2").Select()
